# Append the newest log timestamp to the next empty row in column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
$nextRow = $lastRow + 1

$ws.Cells.Item($nextRow, 1).Value = "2025-10-14 21:10:42"
